# cap nhat quy trinh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account row: same URL as row 2, a new admin login and a new password
$ws.Range("A3").Value = "https://motcua.vinhthanh.cantho.gov.vn/"
$ws.Range("A3").HorizontalAlignment = -4108
$ws.Range("B3").Value = "admin_vinhtrinh"
$ws.Range("C3").Value = "Unitech@"

# Password looks like an email alias, so it gets linked the same way the
# row-2 password is (mailto hyperlink), keeping the cell centered like C2.
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Unitech@")
$ws.Range("C3").HorizontalAlignment = -4108

# Column A needs to be a bit wider to fit the URL text comfortably
$ws.Columns("A:A").ColumnWidth = 35.346354166666664

# Last selection made on the sheet
$ws.Range("D9").Select()
